$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the "Task Progress Daily Estimate" / "Ideal Task Progress" header labels ---
$ws.Range("P3").Value = "Ideal Task Progress"
$ws.Range("AB3").Value = "Task Progress Daily Estimate"

# --- Row 6 (Shan) : Hours worked + Task Progress Daily / Ideal estimate ---
$ws.Range("B6").Value = 1.5
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0.5
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("AD6").Value = 80
$ws.Range("AE6").Value = 80

# --- Row 7 (Pablo) ---
$ws.Range("F7").Value = 2.5
$ws.Range("G7").Value = 2.5
$ws.Range("AD7").Value = 40
$ws.Range("AE7").Value = 60

# --- Row 8 (Jun) ---
$ws.Range("F8").Value = 2.5
$ws.Range("G8").Value = 2.5
$ws.Range("AD8").Value = 0
$ws.Range("AE8").Value = 0

# --- Row 9 (Pedro) ---
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("F9").Value = 3
$ws.Range("AD9").Value = 0
$ws.Range("AE9").Value = 0

# --- Row 10 (Brian) ---
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("F10").Value = 3
$ws.Range("AD10").Value = 0
$ws.Range("AE10").Value = 0

# --- Row 11 (Young) ---
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("X11").Value = 75
$ws.Range("Y11").Value = 100
$ws.Range("AD11").Value = 10
$ws.Range("AE11").Value = 15
# These three used to hold the now-relocated sample values; clear the
# contents but keep their percentage number format.
$ws.Range("AF11").ClearContents()
$ws.Range("AH11").ClearContents()
$ws.Range("AI11").ClearContents()

# --- Row 22 : Meeting Attendance column headers ---
$ws.Range("B22").Value = "Mtg1"
$ws.Range("E22").Value = "Mtg2"
$ws.Range("G22").Value = "Mtg3"

# --- Rows 25-30 : Meeting attendance values ---
$ws.Range("B25").Value = "Abset"
$ws.Range("E25").Value = "Present"
$ws.Range("G25").Value = "Present"

$ws.Range("B26").Value = "Present"
$ws.Range("E26").Value = "Present"
$ws.Range("G26").Value = "Present"

$ws.Range("B27").Value = "Present"
$ws.Range("E27").Value = "Present"
$ws.Range("G27").Value = "Present"

$ws.Range("B28").Value = "Present"
$ws.Range("E28").Value = "Present"
$ws.Range("G28").Value = "Present"

$ws.Range("B29").Value = "Present"
$ws.Range("E29").Value = "Present"
$ws.Range("G29").Value = "Present"

$ws.Range("B30").Value = "Absent"
$ws.Range("E30").Value = "Absent"
$ws.Range("G30").Value = "Present"

# --- Restore the active selection to B11 ---
$ws.Range("B11").Select()
